# Insert a new row above row 181. This pushes the existing rows 181-229
# down to 182-230 (all of their columns move with them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("181:181").Insert()

# After the insert, row 182 holds what used to be row 181's data. Copy
# that whole row's contents into the newly-blank row 181, then change the
# date (column D) of the new row 181 to the new date value.
$ws.Range("A182:R182").Copy()
$ws.Range("A181:R181").PasteSpecial()

$ws.Range("D181").Value = 44551

Write-Output "done"
